$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain formatted as text so values such as
# "62.531.61" or "  -2.92%  " are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.531.61"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "3.371.67"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "572.61"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "125.28"
$ws.Range("E6").Value = "  -7.01%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.366.85"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  -5.66%  "
$ws.Range("D12").Value = "0.375"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "3.958.48"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "3.382.40"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "0.0000169"
$ws.Range("E16").Value = "  -6.75%  "
$ws.Range("D17").Value = "62.635.50"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "24.69"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").Value = "9.12"
$ws.Range("E19").Value = "  -9.40%  "
$ws.Range("D20").Value = "5.56"
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "358.86"
$ws.Range("E22").Value = "  -8.57%  "
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -5.42%  "
$ws.Range("D24").Value = "3.513.59"
$ws.Range("E24").Value = "  -3.63%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "71.32"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").Value = "  -10.65%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "7.00"
$ws.Range("E29").Value = "  -6.15%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "1.41"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -7.49%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.407.16"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  -6.61%  "
$ws.Range("D36").Value = "22.56"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "5.33"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "166.77"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "6.62"
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  -5.59%  "
$ws.Range("D41").Value = "0.0753"
$ws.Range("E41").Value = "  -5.04%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "41.70"
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").Value = "0.762"
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").Value = "  -5.79%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "1.09"
$ws.Range("E46").Value = "  -7.93%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.53"
$ws.Range("E47").Value = "  -8.33%  "
$ws.Range("D48").Value = "22.27"
$ws.Range("E48").Value = "  -10.74%  "
$ws.Range("D49").Value = "6.56"
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").Value = "2.223.91"
$ws.Range("E50").Value = "  -7.08%  "
$ws.Range("D51").Value = "0.835"
$ws.Range("E51").Value = "  -9.84%  "
